# Auto-update citations & altmetric
# Columns are re-mapped: A=doi (was B), B=pubid (was A), C=altid (unchanged),
# D=AltmetricScore numeric value (was title text). Two new rows (20, 21)
# are appended for newly tracked publications that only have a pubid and
# a zero AltmetricScore so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row.
$ws.Cells.Item(1, 1).Value = "doi"
$ws.Cells.Item(1, 2).Value = "pubid"
$ws.Cells.Item(1, 3).Value = "altid"
$ws.Cells.Item(1, 4).Value = "AltmetricScore"

$rows = @(
    @{Row=2; A="10.26028/cybium/2017-413-003"; B="u-x6o8ySG0sC"; C=""; D=0},
    @{Row=3; A="10.3390/jmse6010024"; B="u5HHmVD_uO8C"; C="90861092"; D=1},
    @{Row=4; A="10.1007/s00338-020-01916-8"; B="9yKSN-GCB0IC"; C="77314046"; D=13},
    @{Row=5; A="10.1371/journal.pbio.3000702"; B="qjMakFHDy7sC"; C="96796121"; D=47},
    @{Row=6; A="10.1038/s41597-020-00711-y"; B="2osOgNQ5qMEC"; C="93205299"; D=3},
    @{Row=7; A="10.1111/gcb.15610"; B="UeHWp8X0CEIC"; C="102574363"; D=18},
    @{Row=8; A="10.1002/ece3.8613"; B="IjCSPb-OGe4C"; C="125050022"; D=8},
    @{Row=9; A="10.1111/ddi.13549"; B="zYLM7Y9cAGgC"; C="129013899"; D=38},
    @{Row=10; A="10.1016/j.scitotenv.2022.157049"; B="Tyk-4Ss8FVUC"; C="130854736"; D=31},
    @{Row=11; A="10.1016/j.ecolind.2023.109952"; B="YsMSGLbcyi4C"; C="142625778"; D=6},
    @{Row=12; A="10.1038/s43247-023-00766-w"; B="eQOLeE2rZwMC"; C="145131808"; D=94},
    @{Row=13; A="10.1038/s41598-023-28945-x"; B="Y0pCki6q_DkC"; C="141934419"; D=41},
    @{Row=14; A="10.1111/gcb.17105"; B="_FxGoFyzp5QC"; C="158391407"; D=37},
    @{Row=15; A="10.1016/j.ecss.2024.108734"; B="LkGwnXOMwfcC"; C="170240233"; D=3},
    @{Row=16; A="10.1007/s00338-024-02544-2"; B="UebtZRa9Y70C"; C="166701640"; D=11},
    @{Row=17; A="10.1111/geb.13926"; B="0EnyYjriUFMC"; C="169435683"; D=13},
    @{Row=18; A="10.1038/s41467-025-55949-0"; B="5nxA0vEk-isC"; C="174135131"; D=99},
    @{Row=19; A="10.1007/s10021-025-00995-4"; B="4TOpqqG69KYC"; C="180489954"; D=5},
    @{Row=20; A=""; B="ULOm3_A8WrAC"; C=""; D=0},
    @{Row=21; A=""; B="Zph67rFs4hoC"; C=""; D=0}
)

# Column C (altid) is untouched by this edit -- every existing value stays
# exactly as-is (and rows without an altid still have none), so we only
# ever write columns A, B and D here.
foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($r.A -ne "") {
        $ws.Cells.Item($rowNum, 1).Value = $r.A
    } else {
        $ws.Cells.Item($rowNum, 1).Value = $null
    }

    $ws.Cells.Item($rowNum, 2).Value = $r.B

    $ws.Cells.Item($rowNum, 4).Value = $r.D
}
